# Apply cryptos list update (price + volume refresh, and row reorder for rows 42-45)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'58.819.21"
$ws.Range('E2').Value = '  -3.37%  '
$ws.Range('D3').Value = "'2.557.11"
$ws.Range('E3').Value = '  -1.73%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = "'506.41"
$ws.Range('E5').Value = '  -3.30%  '
$ws.Range('D6').Value = "'142.96"
$ws.Range('E6').Value = '  -7.74%  '
$ws.Range('D7').Value = "'0.998"
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = "'0.552"
$ws.Range('E8').Value = '  -6.22%  '
$ws.Range('D9').Value = "'2.557.49"
$ws.Range('E9').Value = '  -1.93%  '
$ws.Range('D10').Value = "'6.20"
$ws.Range('E10').Value = '  -7.42%  '
$ws.Range('D11').Value = "'0.101"
$ws.Range('E11').Value = '  -4.00%  '
$ws.Range('D12').Value = "'0.329"
$ws.Range('E12').Value = '  -5.19%  '
$ws.Range('E13').Value = '  -0.94%  '
$ws.Range('D14').Value = "'3.002.25"
$ws.Range('E14').Value = '  -1.86%  '
$ws.Range('D15').Value = "'58.817.40"
$ws.Range('E15').Value = '  -3.43%  '
$ws.Range('D16').Value = "'20.57"
$ws.Range('E16').Value = '  -5.07%  '
$ws.Range('D17').Value = "'0.0000134"
$ws.Range('E17').Value = '  -5.02%  '
$ws.Range('D18').Value = "'2.553.87"
$ws.Range('E18').Value = '  -2.08%  '
$ws.Range('D19').Value = "'4.51"
$ws.Range('E19').Value = '  -5.44%  '
$ws.Range('D20').Value = "'330.80"
$ws.Range('E20').Value = '  -6.92%  '
$ws.Range('D21').Value = "'10.04"
$ws.Range('E21').Value = '  -5.19%  '
$ws.Range('D22').Value = "'1.00"
$ws.Range('E22').Value = '  +0.25%  '
$ws.Range('D23').Value = "'5.92"
$ws.Range('E23').Value = '  -4.60%  '
$ws.Range('D24').Value = "'59.31"
$ws.Range('E24').Value = '  -2.82%  '
$ws.Range('D25').Value = "'0.405"
$ws.Range('E25').Value = '  -4.81%  '
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('E27').Value = '  -6.29%  '
$ws.Range('D28').Value = "'0.0₃0774"
$ws.Range('E28').Value = '  -8.36%  '
$ws.Range('D29').Value = "'6.85"
$ws.Range('E29').Value = '  -7.48%  '
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('D31').Value = "'149.38"
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('D32').Value = "'18.54"
$ws.Range('E32').Value = '  -4.52%  '
$ws.Range('D33').Value = "'5.80"
$ws.Range('E33').Value = '  -7.42%  '
$ws.Range('D34').Value = "'1.54"
$ws.Range('E34').Value = '  -3.90%  '
$ws.Range('E35').Value = '  -6.73%  '
$ws.Range('D36').Value = "'0.882"
$ws.Range('E36').Value = '  -3.74%  '
$ws.Range('E37').Value = '  -8.36%  '
$ws.Range('D38').Value = "'35.80"
$ws.Range('E38').Value = '  -1.87%  '
$ws.Range('D39').Value = "'0.824"
$ws.Range('E39').Value = '  -8.59%  '
$ws.Range('D40').Value = "'285.84"
$ws.Range('E40').Value = '  -2.04%  '
$ws.Range('D41').Value = "'1.38"
$ws.Range('E41').Value = '  -8.16%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = "'3.48"
$ws.Range('E42').Value = '  -8.54%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = "'0.997"
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value = "'0.605"
$ws.Range('E44').Value = '  -2.92%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').Value = "'0.0977"
$ws.Range('E45').Value = '  -3.36%  '
$ws.Range('D46').Value = "'0.0529"
$ws.Range('E46').Value = '  -5.54%  '
$ws.Range('D47').Value = "'10.33"
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('D48').Value = "'18.57"
$ws.Range('E48').Value = '  -5.45%  '
$ws.Range('D49').Value = "'0.0226"
$ws.Range('E49').Value = '  -4.98%  '
$ws.Range('D50').Value = "'4.52"
$ws.Range('E50').Value = '  -8.17%  '
$ws.Range('D51').Value = "'1.906.44"
$ws.Range('E51').Value = '  -2.72%  '
